$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.436.78"
$ws.Range("E2").Value = "  +7.41%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.585.86"

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "415.51"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.19"
$ws.Range("E6").Value = "  -0.51%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.651"
$ws.Range("E7").Value = "  +3.68%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.579.34"
$ws.Range("E8").Value = "  +3.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.998"
$ws.Range("E9").Value = "  -0.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.769"
$ws.Range("E10").Value = "  +6.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.176"
$ws.Range("E11").Value = "  +14.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000330"
$ws.Range("E12").Value = "  +43.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.20"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.87"
$ws.Range("E14").Value = "  +1.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.155.71"
$ws.Range("E15").Value = "  +3.47%  "

$ws.Range("E16").Value = "  -0.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.41"
$ws.Range("E17").Value = "  -0.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.586.92"

$ws.Range("E19").Value = "  +5.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.238.68"
$ws.Range("E20").Value = "  +7.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.24"
$ws.Range("E21").Value = "  -2.77%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "454.54"
$ws.Range("E22").Value = "  -1.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "88.59"
$ws.Range("E23").Value = "  -1.92%  "

$ws.Range("E24").Value = "  -4.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.20"
$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("E26").Value = "  +0.90%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.02"
$ws.Range("E27").Value = "  -6.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "34.78"
$ws.Range("E28").Value = "  +4.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.87"
$ws.Range("E29").Value = "  +1.65%  "

$ws.Range("E30").Value = "  +4.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.31"
$ws.Range("E31").Value = "  +1.94%  "

$ws.Range("E32").Value = "  +4.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.37"
$ws.Range("E33").Value = "  -2.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.162"
$ws.Range("E34").Value = "  -4.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.08"
$ws.Range("E35").Value = "  +0.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.77"
$ws.Range("E37").Value = "  -2.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0493"
$ws.Range("E38").Value = "  +0.77%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0746"
$ws.Range("E39").Value = "  +32.73%  "

$ws.Range("E40").Value = "  +9.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.05"
$ws.Range("E42").Value = "  -0.57%  "

$ws.Range("E43").Value = "  +1.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "149.33"
$ws.Range("E44").Value = "  -0.33%  "

$ws.Range("E45").Value = "  -2.42%  "

$ws.Range("E46").Value = "  -2.17%  "

$ws.Range("E47").Value = "  -2.47%  "

$ws.Range("E48").Value = "  -4.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.33"
$ws.Range("E49").Value = "  -2.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.62"
$ws.Range("E50").Value = "  +12.18%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "114.89"
$ws.Range("E51").Value = "  +5.44%  "
